# Fix typo in presentation
#
# 1) "github.com/AnthonyLam/solaaaar" -> "github.com/AnthonyLam/solaaar"
#    (Slide 2, Shape "Shape 49", second paragraph)
# 2) Table style id swaps on the tables in slides 13, 14 (x2) and 15.

$p = $ppt.ActivePresentation

# --- Fix the typo'd URL text -------------------------------------------
$slide2 = $p.Slides.Item(2)
$shape49 = $slide2.Shapes.Item(3)
$tr = $shape49.TextFrame.TextRange
$found = $tr.Find("github.com/AnthonyLam/solaaaar")
$found.Text = "github.com/AnthonyLam/solaaar"

# --- Swap the table style GUIDs ----------------------------------------
$slide13 = $p.Slides.Item(13)
$slide13.Shapes.Item(2).Table.ApplyStyle("{E550B581-EEEF-4063-A7EF-232522EA4822}")

$slide14 = $p.Slides.Item(14)
$slide14.Shapes.Item(2).Table.ApplyStyle("{4419A2C0-D835-4BCE-A1FB-1C9BD85453E1}")
$slide14.Shapes.Item(3).Table.ApplyStyle("{69738B2E-7E7D-43BE-8C86-7E09BE702A0F}")

$slide15 = $p.Slides.Item(15)
$slide15.Shapes.Item(2).Table.ApplyStyle("{ABB37D37-B2F5-4643-AA64-6DAA8B66ADBD}")
